# The "Embodied Carbon" sheet (4th tab) gained a new "DRAM (GPU)" row that
# splits out from the previous single "DRAM" row, which becomes "DRAM (CPU)".
# This shifts every row from the old row 9 downward by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Insert a new blank row above the old row 9 (DRAM); this shifts the old
# rows 9-20 down to 10-21 and keeps the SUM() formula's range in sync.
$ws.Rows.Item(9).Insert()

# The freshly inserted row 9 picked up formatting/placeholder cells from the
# row above it - strip those back to a clean, unstyled row before writing
# the new label into column A only.
$ws.Range("A9:N9").ClearFormats()
$ws.Range("A9:N9").ClearContents()

# Write the new label for the row that used to hold the DRAM (CPU/Samsung)
# data first, then the brand-new "DRAM (GPU)" label - this order matches the
# shared-string table order used in the published workbook.
$ws.Cells.Item(10, 1).Value = "DRAM (CPU)"
$ws.Cells.Item(9, 1).Value = "DRAM (GPU)"

# The picture anchored over columns P:U used to float beside rows 10-16
# (0-indexed 9-15); now that a row was inserted above it, it should float
# beside rows 11-17 (0-indexed 10-16) - nudge it down by one row height
# while preserving its on-sheet size.
$pic = $ws.Shapes.Item(3)
$pic.Top = 134.8
$pic.Left = 822
$pic.Width = 247.2
$pic.Height = 76.6511811023622

# Restore the cursor position that was active when the file was last saved.
$ws.Activate()
$ws.Range("H14").Select()
